{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the async (context) => { ... } function content.\n\nconst body = context.document.body;\n\n// --- 1) Three small in-place text replacements -----------------------\n// Use Range.search (exact match) so we only touch the run that holds the\n// tuple text, preserving its formatting (font/rtl run properties, etc.).\nconst replacements = [\n  {\n    find: \"(Resource, Resource, Resource, Resource);\",\n    replace: \"(Kinds, Resource, Resource, Resource);\"\n  },\n  {\n    find: \"(Kind, Kind, Kind, Kind);\",\n    replace: \"(Statement, Kind, Kind, Kind);\"\n  },\n  {\n    find: \"(Context, Statement, Context, Transform);\",\n    replace: \"(Context, Transform, Transform, Transform);\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n  results.items[0].insertText(replace, \"Replace\");\n}\nawait context.sync();\n\n// --- 2) Insert a new block of paragraphs after the paragraph that ----\n//        contains \"(Mapping, Resource, ResourceMember / Op, Value);\"\nconst anchorResults = body.search(\n  \"(Mapping, Resource, ResourceMember / Op, Value);\",\n  { matchCase: true, matchWholeWord: false }\n);\nanchorResults.load(\"items\");\nawait context.sync();\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found.\");\n}\n\n// New paragraph texts, in document order. Empty strings become blank\n// paragraphs, matching the diff exactly.\nconst newParagraphs = [\n  \"\",\n  \"Augmentations:\",\n  \"Contexts matching Statements applied to aggregated Mapping Context Transforms. \",\n  \"Apply Mappings Transforms. Transform Values Statement (Transform interface reifies Value as Statement Resource)\",\n  \"\",\n  \"Implement Functional APIs:\",\n  \"\",\n  \"Activation (Data)\",\n  \"\",\n  \"Aggregation (Schema)\",\n  \"\",\n  \"Alignment (Behavior)\"\n];\n\n// Insert each paragraph right after the anchor paragraph, one at a time,\n// each subsequent one going after the one just inserted, to preserve order.\nlet anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\nfor (const text of newParagraphs) {\n  anchorParagraph = anchorParagraph.insertParagraph(text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument (aliased $doc / $app) is available.\n\n$d = $word.ActiveDocument\n\n# --- 1) Three small in-place text replacements ------------------------\n# wdReplaceAll = 2 ; wdFindContinue = 1 (wrap param) - use Find/Replace so\n# formatting of the existing run is preserved.\n$replacements = @(\n    @{ Find = \"(Resource, Resource, Resource, Resource);\"; Replace = \"(Kinds, Resource, Resource, Resource);\" },\n    @{ Find = \"(Kind, Kind, Kind, Kind);\"; Replace = \"(Statement, Kind, Kind, Kind);\" },\n    @{ Find = \"(Context, Statement, Context, Transform);\"; Replace = \"(Context, Transform, Transform, Transform);\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n\n# --- 2) Insert a new block of paragraphs after the paragraph that ----\n#        contains \"(Mapping, Resource, ResourceMember / Op, Value);\"\n$anchorRange = $d.Content\n$anchorRange.Find.ClearFormatting()\n$found = $anchorRange.Find.Execute(\"(Mapping, Resource, ResourceMember / Op, Value);\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$anchorParagraph = $anchorRange.Paragraphs(1)\n$anchorIndex = $anchorParagraph.Range.Information(3)\n\n# New paragraph texts, in document order. Empty strings become blank\n# paragraphs, matching the diff exactly.\n$newParagraphs = @(\n    \"\",\n    \"Augmentations:\",\n    \"Contexts matching Statements applied to aggregated Mapping Context Transforms. \",\n    \"Apply Mappings Transforms. Transform Values Statement (Transform interface reifies Value as Statement Resource)\",\n    \"\",\n    \"Implement Functional APIs:\",\n    \"\",\n    \"Activation (Data)\",\n    \"\",\n    \"Aggregation (Schema)\",\n    \"\",\n    \"Alignment (Behavior)\"\n)\n\n$idx = $d.Paragraphs.Count\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq $anchorParagraph.Range.Text) {\n        $idx = $i\n        break\n    }\n}\n\nforeach ($t in $newParagraphs) {\n    $p = $d.Paragraphs($idx)\n    $p.Range.InsertParagraphAfter()\n    $idx = $idx + 1\n    if ($t -ne \"\") {\n        $d.Paragraphs($idx).Range.Text = $t\n    }\n}\n"}
